# Applies the crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.492.28"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "3.469.68"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.95"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.37"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +4.40%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.468.47"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.419"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "4.070.38"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "66.375.14"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "3.467.99"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.54"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.69"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.39"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.537"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  +4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.71"
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.28"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.26"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.887"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.27"
$ws.Range("E39").Value = "  -5.55%  "
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "2.790.40"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.47"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.28"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.30"
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.08"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0292"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("E51").Value = "  -2.15%  "
